# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    holding the per-fund holding breakdown for the new quarter.
# 2. Insert a new leading data row into "总计" summarising that quarter
#    (date/count/value), shifting the previously-existing rows down and
#    renumbering the helper index column.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet, positioned after "2021-Q3" (i.e. right
#    before "总计", which is the last sheet).
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $q3Sheet)
$newSheet.Name = "2022-Q1"

# Match the page margins used by the other quarterly data sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row (row 1) - same column layout/style as the other quarterly sheets.
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Copy the header formatting (bold font + border + centred alignment) from
# the equivalent header row on the "2021-Q3" sheet.
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Data rows (row 2-4): one row per fund holding the stock this quarter.
$fundRows = @(
    @(0, "164811", "工银瑞信中证京津冀协同发展主题指数（LOF）A", 0.23, 94.28, 2.99, 0.0069, 7),
    @(1, "512780", "广发中证京津冀协同发展主题ETF",               0.13, 98.52, 3.13, 0.0041, 8),
    @(2, "164825", "工银瑞信中证京津冀协同发展主题指数（LOF）C", 0.06, 94.28, 2.99, 0.0018, 7)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Column A on the data rows carries the same style as the other sheets'
# index column (e.g. "2021-Q3" A2).
$q3Sheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 2. Update the "总计" (Total) sheet: insert a new top data row for
#    2022-Q1 and shift the existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Inserting a row copies the header row's formatting down onto it; clear
# that back off so the new data row starts unstyled, like the rows below.
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.01

# Column A keeps the same style as the other rows' index column.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)

# Renumber the helper index column (A) for the rows that followed.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3

# Restore the originally active tab ("2021-Q1") - adding/renaming sheets
# above moved the selection onto the new sheet.
$wb.Worksheets.Item("2021-Q1").Activate()
